# Applies the "generic bugfix" edit:
#  - sheet2 ("Ultimo") gains 3 new price-lookup rows for materials that have
#    no recorded cost (S2712T, S1900M, S1103T), each marked "No cost" / "Non trovata" / 0
#  - sheet1 ("log") recalculates the affected totals/weights and messages
#    for the two lavorazioni rows that reference those materials

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet2: insert the 3 new rows (bottom-up so row numbers above stay stable)
# ---------------------------------------------------------------------------

function Insert-PriceRow {
    param($RowIndex, $ColA, $ColB, $ColD)
    $ws2.Rows($RowIndex).Insert() | Out-Null

    # Copy formatting (style s="2") from the row right below the new one so the
    # inserted cells match the rest of the table.
    $ws2.Range("A" + ($RowIndex + 1) + ":E" + ($RowIndex + 1)).Copy() | Out-Null
    $ws2.Range("A" + $RowIndex + ":E" + $RowIndex).PasteSpecial(-4122) | Out-Null

    # Column C always holds the text date "20190103" -- copy it (with its text
    # type preserved) from row 1, which always contains that exact value.
    $ws2.Range("C1").Copy($ws2.Range("C" + $RowIndex))

    $ws2.Range("A" + $RowIndex).Value = $ColA
    $ws2.Range("B" + $RowIndex).Value = $ColB
    $ws2.Range("D" + $RowIndex).Value = $ColD
    $ws2.Range("E" + $RowIndex).Value = 0
}

# Before original row 16 ("Use history" / B6003 / 20190104) -> S1103T
Insert-PriceRow 16 "No cost" "S1103T" "Non trovata"

# Before original row 7 ("Use history" / B6003 / 20190105) -> S1900M
Insert-PriceRow 7 "No cost" "S1900M" "Non trovata"

# Before original row 5 (A0601) -> S2712T
Insert-PriceRow 5 "No cost" "S2712T" "Non trovata"

# ---------------------------------------------------------------------------
# Sheet1: update the two lavorazione summary rows
# ---------------------------------------------------------------------------

$detail2 = @"
Lavorazioni toccate:
 [LAH/004045 q.: 4050.0]
Totale carichi: 4050.0

Lavorazione Linea 5: euro/kg. 0.264 x 4050.0 = 1069.2

Costi materie prime:
Lavoration LAH/004045:
 - A1224: EUR 0.77 x q. 2800.0 = 2156.0
 - A0102: EUR 0.62 x q. 240.0 = 148.8
 - A0404: EUR 0.92 x q. 420.0 = 386.4
 - A0405: EUR 0.325 x q. 420.0 = 136.5
 - S2712T: EUR 0.0 x q. 80.0 = 0.0
 - A0601: EUR 1.09 x q. 32.0 = 34.88
 - A2035: EUR 4.4 x q. 8.0 = 35.2
 - S1900M: EUR 0.0 x q. 400.0 = 0.0
Totale materie prime: 2897.78

Costi imballi e pallet:
 - Imballo [LAH/004045] B6003: EUR 0.405 x q. 162 = 65.61
 - Pallet [LAH/004045] B6107: EUR 10.7 x q. 4 = 42.8 
Totale imballi: 3006.19

Peso materie prime: 4400.0

Costo totale:
EUR 4075.39 : q. 4050.0 = EUR/unit 1.0062691358 (carico)

"@

$ws1.Range("F2").Value = 4400
$ws1.Range("I2").Value = $detail2
$ws1.Range("N2").Value = "Material with price 0, Material with price 0"

$detail3 = @"
Lavorazioni toccate:
 [LAH/004047 q.: 6050.0]
Totale carichi: 6050.0

Lavorazione Linea 4: euro/kg. 0.264 x 6050.0 = 1597.2

Costi materie prime:
Lavoration LAH/004047:
 - A1004: EUR 0.645 x q. 969.0 = 625.005
 - A1216: EUR 0.68 x q. 969.0 = 658.92
 - A0600: EUR 0.98 x q. 229.5 = 224.91
 - A3004: EUR 2.4 x q. 102.0 = 244.8
 - A0402: EUR 0.88 x q. 17.85 = 15.708
 - A0404: EUR 0.92 x q. 384.03 = 353.3076
 - A0403: EUR 0.31 x q. 2428.62 = 752.8722
 - S1103T: EUR 0.0 x q. 925.0 = 0.0
Totale materie prime: 2875.5228

Costi imballi e pallet:
 - Imballo [LAH/004047] B6003: EUR 0.405 x q. 242 = 98.01
 - Pallet [LAH/004047] B6107: EUR 10.7 x q. 6 = 64.2 
Totale imballi: 3037.7328

Peso materie prime: 6025.0

Costo totale:
EUR 4634.9328 : q. 6050.0 = EUR/unit 0.766104595041 (carico)

"@

$ws1.Range("F3").Value = 6025
$ws1.Range("G3").ClearContents()
$ws1.Range("I3").Value = $detail3
$ws1.Range("M3").Value = "X"
$ws1.Range("N3").Value = "Material with price 0"

Write-Host "Edit applied."
